$wb = $excel.ActiveWorkbook

# --- Rename sheets to add quality-profile descriptors ---
$wsUltra  = $wb.Worksheets.Item(1)   # "0.05"  -> "Ultra 0.05"
$wsHigh   = $wb.Worksheets.Item(2)   # "0.1"   -> "High 0.1"
$wsNormal = $wb.Worksheets.Item(3)   # "0.2"   -> "Normal 0.2"
$wsFast   = $wb.Worksheets.Item(4)   # "0.3"   -> "Fast 0.3"
$wsDraft  = $wb.Worksheets.Item(5)   # "0.34"  -> "Draft 0.34"

$wsUltra.Name  = "Ultra 0.05"
$wsHigh.Name   = "High 0.1"
$wsNormal.Name = "Normal 0.2"
$wsFast.Name   = "Fast 0.3"
$wsDraft.Name  = "Draft 0.34"

# --- "Ultra 0.05" sheet: new column B for the (unavailable) Nylon quality profile ---
$wsUltra.Columns.Item(2).ColumnWidth = 51.33

# Highlight C8 (PLA value) the same way the other mismatched rows are highlighted
$wsUltra.Range("C8").Interior.Pattern = -4124
$wsUltra.Range("C17").Copy()
$wsUltra.Range("C8").PasteSpecial(-4122)

# --- Selections / active tab bookkeeping ---
$wsUltra.Range("B12").Select()

$wsFast.Activate()
$wsFast.Range("A49").Select()

$wsDraft.Activate()
$wsDraft.Range("A43").Select()
